$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the old row 30 ("Ograniczenia dolne" header),
# shifting the old rows 30-39 down to 33-42. Each new row starts as a
# copy of row 29 (the last "CAP_BND / ELE_NEW_PV_GRND" row) so it
# inherits the same formatting (including the lone styled I-column cell).
$ws.Rows("29:29").Copy()
$ws.Rows("30:30").Insert()
$ws.Rows("29:29").Copy()
$ws.Rows("30:30").Insert()
$ws.Rows("29:29").Copy()
$ws.Rows("30:30").Insert()

$ws.Rows("30:32").RowHeight = 18.75

# New row 30: UP / CAP_BND / 2030 / 6 / ELE_NEW_PV_GRND
$ws.Range("B30").Value = "UP"
$ws.Range("C30").Value = "CAP_BND"
$ws.Range("D30").Value = 2030
$ws.Range("E30").Value = 6
$ws.Range("F30").Value = "ELE_NEW_PV_GRND"

# New row 31: UP / CAP_BND / 2035 / 10 / ELE_NEW_PV_GRND
$ws.Range("B31").Value = "UP"
$ws.Range("C31").Value = "CAP_BND"
$ws.Range("D31").Value = 2035
$ws.Range("E31").Value = 10
$ws.Range("F31").Value = "ELE_NEW_PV_GRND"

# New row 32: UP / CAP_BND / 2040 / 12 / ELE_NEW_PV_GRND
$ws.Range("B32").Value = "UP"
$ws.Range("C32").Value = "CAP_BND"
$ws.Range("D32").Value = 2040
$ws.Range("E32").Value = 12
$ws.Range("F32").Value = "ELE_NEW_PV_GRND"

# The header row that used to be row 30 ("\I: " / "Ograniczenia dolne
# (wymuszenia)") is now row 33; its height becomes 18 (no longer inherits
# the taller 18.75 used for the data rows around it).
$ws.Rows("33:33").RowHeight = 18

# Reposition the view roughly like the source: scrolled so row 16 is at
# the top, with F31 selected.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("F31").Select()
